$d = $word.ActiveDocument

# The paragraph between "Paragraph start" and "Paragraph end" currently
# spells out "${values}" split across three runs ("${value", "s", "}").
# Collapse it back down to the single-run text "${value}".
$d.Content.Find.Execute('${values}', $false, $false, $false, $false, $false,
                         $true, 1, $false, '${value}', 2)
